$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 399, shifting existing rows 399..466 down to 400..467
$ws.Rows.Item(399).Insert()

# Populate the newly inserted row 399 with the new record
$ws.Cells.Item(399, 1).Value = 4
$ws.Cells.Item(399, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(399, 3).Value = "Los Lagos"
$ws.Cells.Item(399, 4).Value = 45258
$ws.Cells.Item(399, 5).Value = 10
$ws.Cells.Item(399, 6).Value = 100112028
$ws.Cells.Item(399, 7).Value = "Sandia"
$ws.Cells.Item(399, 8).Value = "Sin especificar"
$ws.Cells.Item(399, 9).Value = "Primera"
$ws.Cells.Item(399, 10).Value = 800
$ws.Cells.Item(399, 11).Value = 1100
$ws.Cells.Item(399, 12).Value = 1200
$ws.Cells.Item(399, 13).Value = 1150
$ws.Cells.Item(399, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(399, 15).Value = "Perú"
$ws.Cells.Item(399, 16).Value = 1150
$ws.Cells.Item(399, 17).Value = 1
$ws.Cells.Item(399, 18).Value = "Hortaliza"
